$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('to live', '住む|すむ'),
    @('address', '住所|じゅうしょ'),
    @('to immigrate', '移住する|いじゅうする'),
    @('New Year', 'お正月|おしょうがつ'),
    @('right', '正しい|ただしい'),
    @('noon', '正午|しょうご'),
    @('correct answer', '正解|せいかい'),
    @('third-year student', '三年生|さんねんせい'),
    @('next year', '来年|らいねん'),
    @('this year', '今年|ことし'),
    @('year', '年|とし'),
    @('to sell', '売る|うる'),
    @('stand; stall', '売店|ばいてん'),
    @('vending machine', '自動販売機|じどうはんばいき'),
    @('to buy', '買う|かう'),
    @('shopping', '買い物|かいもの'),
    @('selling and buying', '売買|ばいばい'),
    @('town', '町|まち'),
    @('Kitayama town', '北山町|きたやまちょう'),
    @('mayor of a town', '町長|ちょうちょう'),
    @('long', '長い|ながい'),
    @('the eldest son', '長男|ちょうなん'),
    @('company president', '社長|しゃちょう'),
    @('way; road', '道|みち'),
    @('calligraphy', '書道|しょどう'),
    @('judo', '柔道|じゅうどう'),
    @('Hokkaido', '北海道|ほっかいどう'),
    @('snow', '雪|ゆき'),
    @('new snow', '新雪|しんせつ'),
    @('snowman', '雪だるま|ゆきだるま'),
    @('to stand', '立つ|たつ'),
    @('national university', '国立大学|こくりつだいがく'),
    @('private high school', '私立高校|しりつこうこう'),
    @('oneself', '自分|じぶん'),
    @('automobile', '自動車|じどうしゃ'),
    @('bicycle', '自転車|じてんしゃ'),
    @('freedom', '自由|じゆう'),
    @('night', '夜|よる'),
    @('middle of night; midnight', '夜中|よなか'),
    @('tonight', '今夜|こんや'),
    @('dawn', '夜明け|よあけ'),
    @('morning', '朝|あさ'),
    @('this morning', '今朝|けさ'),
    @('breakfast', '朝食|ちょうしょく'),
    @('every morning', '毎朝|まいあさ'),
    @('to carry; to hold', '持つ|もつ'),
    @('to bring', '持ってくる|もってくる'),
    @('belongings', '所持品|しょじひん'),
    @('feeling', '気持ち|きもち'),
)

$startRow = 93
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
